$d = $word.ActiveDocument

function Find-ParaIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# 1) Remove the "Meta description: ..." paragraph that immediately follows
#    the document title paragraph.
$metaIndex = Find-ParaIndex("Meta description")
if ($metaIndex -ge 1) {
    $d.Paragraphs.Item($metaIndex).Range.Delete()
}

# 2) Insert a new bold paragraph ("Play Asgard Slot Game for Free - Review")
#    right before the final paragraph in the document (the one that used to
#    hold the image-generation prompt text).
$promptIndex = Find-ParaIndex("Can you create an image for Asgard")
$promptPara = $d.Paragraphs.Item($promptIndex)
$insertPoint = $promptPara.Range.Start
$target = $d.Range($insertPoint, $insertPoint)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Asgard Slot Game for Free - Review</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)

# The insertion above splices in the new paragraph together with an extra
# blank spacer paragraph (needed to force a paragraph break); remove the
# spacer, which is the lone empty paragraph left right after the new one.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r") {
        $p.Range.Delete()
        break
    }
}

# 3) Replace the old image-generation prompt text (now in the final
#    paragraph) with the meta-description sentence, preserving the italic
#    run formatting already present on that text.
$find = $d.Content.Find
$find.Execute("Can you create an image for Asgard with the following specifications: - A cartoon style - A happy Maya warrior with glasses as the main character", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  "Discover Asgard Slot Game for Free. Our review covers gameplay mechanics, special features, graphics, and playing modes.", 2)
